# Update database credentials and dependencies
# (Refreshes the lot/roll tracking table with the latest RAPID export.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("A2").Value = "FWFT00073 00006"
$ws.Range("B2").Value = 2305300406
$ws.Range("C2").Value = "W2305300406-01"
$ws.Range("D2").Value = "HK"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 82.11
$ws.Range("G2").Value = 159

# ---- Row 3 ----
$ws.Range("A3").Value = "FWFT00073 00006"
$ws.Range("B3").Value = 2305300406
$ws.Range("C3").Value = "W2305300406-01"
$ws.Range("D3").Value = "HK"
$ws.Range("E3").Value = "4"
$ws.Range("F3").Value = 71.57
$ws.Range("G3").Value = 159

# ---- Row 4 ----
$ws.Range("A4").Value = "FWFT00073 00006"
$ws.Range("B4").Value = 2305300405
$ws.Range("C4").Value = "W2305300405-01"
$ws.Range("D4").Value = "HK"
$ws.Range("E4").Value = "5E"
$ws.Range("F4").Value = 12.22
$ws.Range("G4").Value = 158

# ---- Row 5 ----
$ws.Range("A5").Value = "FWFT00073 00006"
$ws.Range("B5").Value = 2305300405
$ws.Range("C5").Value = "W2305300405-01"
$ws.Range("D5").Value = "HK"
$ws.Range("E5").Value = "14"
$ws.Range("F5").Value = 84.16
$ws.Range("G5").Value = 158

# ---- Row 6 ----
$ws.Range("A6").Value = "FWFT00073 00006"
$ws.Range("B6").Value = 2305300410
$ws.Range("C6").Value = "W2305300410-01"
$ws.Range("D6").Value = "HK"
$ws.Range("E6").Value = "5B"
$ws.Range("F6").Value = 12.6
$ws.Range("G6").Value = 158

# ---- Row 7 (was blank, now a new data row) ----
$ws.Range("A7").Value = "FWFT00073 00006"
$ws.Range("B7").Value = 2305300410
$ws.Range("C7").Value = "W2305300410-01"
$ws.Range("D7").Value = "HK"
$ws.Range("E7").Value = "9"
$ws.Range("F7").Value = 63.64
$ws.Range("G7").Value = 169

# Left-align the Roll No column (E) to match the updated formatting.
$ws.Columns("E").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft

# Move the active selection to where the user left off editing.
$ws.Range("B8").Select()
